# Add gender, age quota AU, CA, TR, UA
# Strategy: duplicate the existing "quotas_SK" template sheet 4 times (it already
# carries the correct styles/layout/formulas for the age-quota block), rename
# each copy, then overwrite the country-specific inputs (gender split C2, age
# quota B8:F8) and switch the gender-ratio formula from the "F3 ratio" pattern
# used by the older quota sheets to the "1-C2" pattern used by these new ones.

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("quotas_SK")

$newSheets = @(
    @{ Name = "quota_AU"; C2 = 0.49354379369142626;  B8 = 0.11192255877789467; C8 = 0.18591237270266067; D8 = 0.26180562458671341; E8 = 0.23014443705193119; F8 = 0.21021500688079983 },
    @{ Name = "quota_CA"; C2 = 0.4927381777223736;   B8 = 0.10402815988780871; C8 = 0.17502163559091188; D8 = 0.24488811098766669; E8 = 0.25292269475436224; F8 = 0.22313939877924996 },
    @{ Name = "quota_TR"; C2 = 0.48657573802133475;  B8 = 0.1582325601298683;  C8 = 0.2131137356790879;  D8 = 0.29681998391047659; E8 = 0.20571250504618752; F8 = 0.12612121523437902 },
    @{ Name = "quota_UA"; C2 = 0.45142595728437557;  B8 = 0.082159805419933327; C8 = 0.17834324173209759; D8 = 0.28227842928850105; E8 = 0.24861654266018471; F8 = 0.20860198089928325 }
)

foreach ($info in $newSheets) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $template.Copy($null, $lastSheet)
    $ws = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws.Name = $info.Name

    # Gender split: new sheets type the female share directly into C2 and
    # derive the male share as 1-C2 (older quota sheets instead derived both
    # from a M/F ratio typed into F3).
    $ws.Range("B2").Formula = "=1-C2"
    $ws.Range("C2").Value = $info.C2
    $ws.Range("F3").ClearContents()

    # Age quota shares.
    $ws.Range("B8").Value = $info.B8
    $ws.Range("C8").Value = $info.C8
    $ws.Range("D8").Value = $info.D8
    $ws.Range("E8").Value = $info.E8
    $ws.Range("F8").Value = $info.F8

    $ws.Range("C2").Select()
}

# View bookkeeping to mirror the author's final click trail.
$spec = $wb.Worksheets.Item("Specificities")
$spec.Activate()
$spec.Range("R94").Select()
$excel.ActiveWindow.ScrollRow = 2
$spec.Range("V4").Select()

$ch = $wb.Worksheets.Item("quotas_CH")
$ch.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ch.Range("F3").Select()

$mx = $wb.Worksheets.Item("quotas_MX")
$mx.Activate()
$mx.Range("I18").Select()

$ua = $wb.Worksheets.Item("quota_UA")
$ua.Activate()
$ua.Range("H4").Select()
